$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 66670640
$ws.Range("I64").Value = 200002670
$ws.Range("J64").Value = 4622
$ws.Range("K64").Value = 200002670
$ws.Range("L64").Value = 4622
$ws.Range("M64").Value = -200002422
$ws.Range("N64").Value = -5118
$ws.Range("H67").Value = 66670640
$ws.Range("I67").Value = 200002670
$ws.Range("J67").Value = 4622
$ws.Range("K67").Value = 200002670
$ws.Range("L67").Value = 4622
$ws.Range("M67").Value = -200001812
$ws.Range("N67").Value = -6338
$ws.Range("H93").Value = 32667.334
$ws.Range("J93").Value = 32667.334
$ws.Range("L93").Value = 32667.334
$ws.Range("N93").Value = -37659.334
$ws.Range("H135").Value = 577.2941
$ws.Range("I135").Value = 515.8461
$ws.Range("J135").Value = 777
$ws.Range("K135").Value = 4642.6149
$ws.Range("L135").Value = 6993
$ws.Range("M135").Value = -2107.6149
$ws.Range("N135").Value = -12063
$ws.Range("H137").Value = 2500.6025
$ws.Range("I137").Value = 882.08826
$ws.Range("J137").Value = 3751.2727
$ws.Range("K137").Value = 2646.26478
$ws.Range("L137").Value = 11253.8181
$ws.Range("M137").Value = -96.26477999999997
$ws.Range("N137").Value = -16353.8181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1768.1177
$ws.Range("I61").Value = 1462.7391
$ws.Range("K61").Value = 1462.7391
$ws.Range("M61").Value = -1250.7391
$ws.Range("H74").Value = 8041.857
$ws.Range("I74").Value = 925.8
$ws.Range("J74").Value = 25832
$ws.Range("K74").Value = 925.8
$ws.Range("L74").Value = 25832
$ws.Range("M74").Value = -51.79999999999995
$ws.Range("N74").Value = -27580
$ws.Range("H77").Value = 8041.857
$ws.Range("I77").Value = 925.8
$ws.Range("J77").Value = 25832
$ws.Range("K77").Value = 4629
$ws.Range("L77").Value = 129160
$ws.Range("M77").Value = -261
$ws.Range("N77").Value = -137896
$ws.Range("H132").Value = 9872.154
$ws.Range("I132").Value = 7525
$ws.Range("J132").Value = 15153.25
$ws.Range("K132").Value = 22575
$ws.Range("L132").Value = 45459.75
$ws.Range("M132").Value = -20045
$ws.Range("N132").Value = -50519.75
$ws.Range("H136").Value = 1768.1177
$ws.Range("I136").Value = 1462.7391
$ws.Range("K136").Value = 4388.2173
$ws.Range("M136").Value = -1838.2173

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 13878
$ws.Range("I102").Value = 2170.6667
$ws.Range("J102").Value = 49000
$ws.Range("K102").Value = 2170.6667
$ws.Range("L102").Value = 49000
$ws.Range("M102").Value = 1074.3333
$ws.Range("N102").Value = -55490
$ws.Range("H134").Value = 1100.2593
$ws.Range("I134").Value = 912.73914
$ws.Range("J134").Value = 2178.5
$ws.Range("K134").Value = 2738.21742
$ws.Range("L134").Value = 6535.5
$ws.Range("M134").Value = -203.2174199999999
$ws.Range("N134").Value = -11605.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20841.207
$ws.Range("I31").Value = 1126.0834
$ws.Range("J31").Value = 34757.766
$ws.Range("K31").Value = 1126.0834
$ws.Range("L31").Value = 34757.766
$ws.Range("M31").Value = -831.0834
$ws.Range("N31").Value = -35347.766
$ws.Range("H34").Value = 20841.207
$ws.Range("I34").Value = 1126.0834
$ws.Range("J34").Value = 34757.766
$ws.Range("K34").Value = 1126.0834
$ws.Range("L34").Value = 34757.766
$ws.Range("M34").Value = -924.0834
$ws.Range("N34").Value = -35161.766
$ws.Range("H58").Value = 1205
$ws.Range("I58").Value = 1083.3889
$ws.Range("J58").Value = 1642.8
$ws.Range("K58").Value = 1083.3889
$ws.Range("L58").Value = 1642.8
$ws.Range("M58").Value = -880.3888999999999
$ws.Range("N58").Value = -2048.8
$ws.Range("H92").Value = 22288
$ws.Range("J92").Value = 22288
$ws.Range("L92").Value = 22288
$ws.Range("N92").Value = -27280
$ws.Range("H132").Value = 50008776
$ws.Range("I132").Value = 83344890
$ws.Range("J132").Value = 4612
$ws.Range("K132").Value = 250034670
$ws.Range("L132").Value = 13836
$ws.Range("M132").Value = -250032140
$ws.Range("N132").Value = -18896
$ws.Range("H134").Value = 1829.1072
$ws.Range("I134").Value = 1947.6111
$ws.Range("J134").Value = 1615.8
$ws.Range("K134").Value = 5842.8333
$ws.Range("L134").Value = 4847.4
$ws.Range("M134").Value = -3307.8333
$ws.Range("N134").Value = -9917.4
$ws.Range("H136").Value = 1205
$ws.Range("I136").Value = 1083.3889
$ws.Range("J136").Value = 1642.8
$ws.Range("K136").Value = 3250.1667
$ws.Range("L136").Value = 4928.4
$ws.Range("M136").Value = -700.1666999999998
$ws.Range("N136").Value = -10028.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 303
$ws.Range("I2").Value = 394.45834
$ws.Range("J2").Value = 28.625
$ws.Range("K2").Value = 2366.75004
$ws.Range("L2").Value = 171.75
$ws.Range("M2").Value = -2253.75004
$ws.Range("N2").Value = -397.75
$ws.Range("H3").Value = 3838.5173
$ws.Range("I3").Value = 2014.9166
$ws.Range("J3").Value = 12591.8
$ws.Range("K3").Value = 6044.7498
$ws.Range("L3").Value = 37775.39999999999
$ws.Range("M3").Value = -5932.7498
$ws.Range("N3").Value = -37999.39999999999
$ws.Range("H131").Value = 871.46155
$ws.Range("I131").Value = 439.14285
$ws.Range("J131").Value = 938.7111
$ws.Range("K131").Value = 1317.42855
$ws.Range("L131").Value = 2816.1333
$ws.Range("M131").Value = 3722.57145
$ws.Range("N131").Value = -12896.1333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 34195
$ws.Range("J88").Value = 34195
$ws.Range("L88").Value = 34195
$ws.Range("N88").Value = -35097
$ws.Range("H91").Value = 34195
$ws.Range("J91").Value = 34195
$ws.Range("L91").Value = 34195
$ws.Range("N91").Value = -37315
$ws.Range("H132").Value = 9574.866
$ws.Range("I132").Value = 11052.167
$ws.Range("J132").Value = 3665.6667
$ws.Range("K132").Value = 33156.501
$ws.Range("L132").Value = 10997.0001
$ws.Range("M132").Value = -30626.501
$ws.Range("N132").Value = -16057.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5626.7
$ws.Range("I132").Value = 7704.579
$ws.Range("J132").Value = 2037.6364
$ws.Range("K132").Value = 23113.737
$ws.Range("L132").Value = 6112.9092
$ws.Range("M132").Value = -20583.737
$ws.Range("N132").Value = -11172.9092
$ws.Range("H136").Value = 6292.048
$ws.Range("I136").Value = 1566.3529
$ws.Range("J136").Value = 26376.25
$ws.Range("K136").Value = 4699.0587
$ws.Range("L136").Value = 79128.75
$ws.Range("M136").Value = -2149.0587
$ws.Range("N136").Value = -84228.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7389.5264
$ws.Range("I132").Value = 9135.929
$ws.Range("J132").Value = 2499.6
$ws.Range("K132").Value = 27407.787
$ws.Range("L132").Value = 7498.799999999999
$ws.Range("M132").Value = -24877.787
$ws.Range("N132").Value = -12558.8
$ws.Range("H136").Value = 3809.842
$ws.Range("I136").Value = 5292.476
$ws.Range("J136").Value = 1978.3529
$ws.Range("K136").Value = 15877.428
$ws.Range("L136").Value = 5935.0587
$ws.Range("M136").Value = -16877.428
$ws.Range("N136").Value = -11035.0587
